$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Set all the flow result booleans (Create/Read/Update/Delete Test Passed)
# for every device/zone/category row to TRUE, so that "running the main
# flow" marks every individual flow as having passed in one click.
$ws.Range("B2:E24").Value = $true
